$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header strings (Volume/Number and Report Covering Week dates) ---
$ws.Range("A8").Value = "Volume 29   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# --- Update crime-statistics table (rows 14-30) ---
$ws.Range("C14").Value = 1
$ws.Range("G14").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$ws.Range("F14").Value = 1
$ws.Range("G14").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("H14").Value = -50
$ws.Range("I14").Value = 7
$ws.Range("K14").Value = 40
$ws.Range("L14").Value = 133.333333333333
$ws.Range("N14").Value = -30
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = 17.647058823529
$ws.Range("N15").Value = -48.717948717948
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 61
$ws.Range("G16").Value = 36
$ws.Range("H16").Value = 69.444444444444
$ws.Range("I16").Value = 617
$ws.Range("J16").Value = 410
$ws.Range("K16").Value = 50.487804878048
$ws.Range("L16").Value = 290.506329113924
$ws.Range("M16").Value = 303.267973856209
$ws.Range("N16").Value = -74.140821458508
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = "0"
$ws.Range("A17").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Range("E17").Value = "***.*"
$ws.Range("A17").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = 46.153846153846
$ws.Range("I17").Value = 453
$ws.Range("K17").Value = 7.857142857142
$ws.Range("L17").Value = 119.902912621359
$ws.Range("M17").Value = 157.386363636364
$ws.Range("N17").Value = -28.885400313971
$ws.Range("C18").Value = 16
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = 77.777777777777
$ws.Range("F18").Value = 58
$ws.Range("G18").Value = 46
$ws.Range("H18").Value = 26.086956521739
$ws.Range("I18").Value = 632
$ws.Range("J18").Value = 420
$ws.Range("K18").Value = 50.476190476190
$ws.Range("L18").Value = 86.982248520710
$ws.Range("M18").Value = 92.09726443769
$ws.Range("N18").Value = -75.166994106090
$ws.Range("C19").Value = 64
$ws.Range("D19").Value = 40
$ws.Range("E19").Value = 60
$ws.Range("F19").Value = 238
$ws.Range("G19").Value = 159
$ws.Range("H19").Value = 49.685534591195
$ws.Range("I19").Value = 2206
$ws.Range("J19").Value = 1346
$ws.Range("K19").Value = 63.893016344725
$ws.Range("L19").Value = 99.818840579710
$ws.Range("M19").Value = 2.461681374825
$ws.Range("N19").Value = -75.210697831217
$ws.Range("C20").Value = "0"
$ws.Range("A17").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 52
$ws.Range("K20").Value = 21.153846153846
$ws.Range("L20").Value = 43.181818181818
$ws.Range("N20").Value = -81.415929203539
$ws.Range("C21").Value = 104
$ws.Range("D21").Value = 56
$ws.Range("E21").Value = 85.714285714285
$ws.Range("F21").Value = 403
$ws.Range("G21").Value = 276
$ws.Range("H21").Value = 46.014492753623
$ws.Range("I21").Value = 3998
$ws.Range("J21").Value = 2670
$ws.Range("K21").Value = 49.737827715355
$ws.Range("L21").Value = 112.999467234949
$ws.Range("M21").Value = 40.280701754386
$ws.Range("N21").Value = -73.086502860989
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -60
$ws.Range("F22").Value = 7
$ws.Range("G22").Value = 18
$ws.Range("H22").Value = -61.111111111111
$ws.Range("I22").Value = 174
$ws.Range("J22").Value = 150
$ws.Range("K22").Value = 16
$ws.Range("L22").Value = 35.9375
$ws.Range("M22").Value = 24.285714285714
$ws.Range("C24").Value = 55
$ws.Range("D24").Value = 56
$ws.Range("E24").Value = -1.785714285714
$ws.Range("F24").Value = 283
$ws.Range("G24").Value = 223
$ws.Range("H24").Value = 26.905829596412
$ws.Range("I24").Value = 3211
$ws.Range("J24").Value = 2152
$ws.Range("K24").Value = 49.210037174721
$ws.Range("L24").Value = 87.558411214953
$ws.Range("M24").Value = -29.875518672199
$ws.Range("C25").Value = 30
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 172.727272727273
$ws.Range("F25").Value = 77
$ws.Range("G25").Value = 57
$ws.Range("H25").Value = 35.087719298245
$ws.Range("I25").Value = 870
$ws.Range("J25").Value = 826
$ws.Range("K25").Value = 5.326876513317
$ws.Range("L25").Value = 74.698795180722
$ws.Range("M25").Value = 52.097902097902
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 34
$ws.Range("K26").Value = -14.705882352941
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 600
$ws.Range("F27").Value = 23
$ws.Range("G27").Value = 14
$ws.Range("H27").Value = 64.285714285714
$ws.Range("I27").Value = 217
$ws.Range("J27").Value = 146
$ws.Range("K27").Value = 48.630136986301
$ws.Range("L27").Value = 92.035398230088
$ws.Range("D30").Value = 1
$ws.Range("G14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = -100
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("J30").Value = 36
$ws.Range("K30").Value = -36.111111111111
